$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after row 765 (new row becomes row 766), shifting existing
# rows 766:807 down to 767:808.
$ws.Rows.Item(766).Insert()

# Populate the newly inserted row 766 with the new data point.
# Column A holds the date as plain text (matching the rest of the column),
# so force text format before assigning to avoid Excel auto-converting the
# string into a date serial number.
$ws.Cells.Item(766, 1).NumberFormat = "@"
$ws.Cells.Item(766, 1).Value = "2026/02/02"
$ws.Cells.Item(766, 1).Style = "Normal"
$ws.Cells.Item(766, 2).Value = "月"
$ws.Cells.Item(766, 3).Value = 6
$ws.Cells.Item(766, 4).Value = 161
